$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (ram / aiautomationhig@gmail.com / Wealth Manager / ...) -- it's being removed,
# shifting the old row 3 up to row 2, and old row 4 up to row 3.
$ws.Rows.Item(2).Delete()

# Update the new last row (formerly row 4, now row 3) with the new values.
# Force text format on numeric-looking fields so they stay stored as text.
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "Abu Inshah"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "9943374466"
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = "aiautomationhig@gmail.com"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "Wealth Manager"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "uploads/abu_inshah_1752329088201.jpeg"
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = "1752329088240"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "2025-07-12T14:04:48.240Z"
